$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D and E columns are treated as text so numeric/percent-looking
# strings are not auto-converted to numbers by the COM layer.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "297.70"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-2.00%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "31.34"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.67%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.123"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-2.70%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07330"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-2.23%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "7.752"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-1.03%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.734"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "17.13%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.728"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.97%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9230"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.91%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1668"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-1.36%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06943"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-8.00%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07977"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.16%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03017"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.88%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09924"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.18%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001489"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.25%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006175"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-3.63%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.456"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.11%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.220"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.51%"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3227"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-2.52%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1330"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.00%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.552"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.70%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04655"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "2.07%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.1581"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-3.00%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001218"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.25%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004738"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "6.58%"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-7.26%"

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "7.73%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01713"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "3.19%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04455"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.89%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007139"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.09%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1333"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-1.14%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002206"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.90%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01079"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-17.08%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006107"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-1.91%"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-21.34%"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "170.80%"
